$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("procedimientos")

$ws.Range("A9").Value = "proceso"
$ws.Range("B9").Value = "excluir_por_morosidad"
$ws.Range("C9").Value = "Procedimiento para calcular si las operaciones morosas tienen 5 años o mas"

$ws.Range("C10").Select()
